$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F23").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F25").Value = "ppe"
$ws.Range("F26").Value = "ppe"
$ws.Range("F29").Value = "application instructions"
$ws.Range("F30").Value = "env warning - species || pollinator"
$ws.Range("F32").Value = "env warning - water"
$ws.Range("F33").Value = "env warning - water || off target movement"
$ws.Range("F34").Value = "env warning - water || off target movement"
$ws.Range("F37").Value = "application instructions"
$ws.Range("F38").Value = "application instructions"
$ws.Range("F39").Value = "application instructions"
$ws.Range("F40").Value = "135_product_information"
$ws.Range("F41").Value = "use restrictions"
$ws.Range("F42").Value = "use restrictions"
$ws.Range("F49").Value = "application instructions"
$ws.Range("F50").Value = "application instructions"
$ws.Range("F51").Value = "application instructions"
$ws.Range("F52").Value = "application instructions"
$ws.Range("F54").Value = "application instructions"
$ws.Range("F55").Value = "application instructions"
$ws.Range("F56").Value = "application instructions"
$ws.Range("F57").Value = "application instructions"
$ws.Range("F63").Value = "application instructions"
$ws.Range("F64").Value = "irrigation"
$ws.Range("F66").Value = "application instructions"
$ws.Range("F67").Value = "application instructions"
$ws.Range("F68").Value = "application instructions"
$ws.Range("F70").Value = "application instructions"
$ws.Range("F71").Value = "application instructions"
$ws.Range("F73").Value = "off target movement"
$ws.Range("F74").Value = "off target movement"
$ws.Range("F75").Value = "off target movement"
$ws.Range("F77").Value = "off target movement"
$ws.Range("F78").Value = "off target movement"
$ws.Range("F80").Value = "off target movement"
$ws.Range("F82").Value = "off target movement"
$ws.Range("F83").Value = "off target movement"
$ws.Range("F84").Value = "off target movement"
$ws.Range("F85").Value = "off target movement"
$ws.Range("F86").Value = "off target movement"
$ws.Range("F87").Value = "off target movement"
$ws.Range("F88").Value = "off target movement"
$ws.Range("F89").Value = "off target movement"
$ws.Range("F90").Value = "mixing"
$ws.Range("F91").Value = "mixing"
$ws.Range("F93").Value = "mixing"
$ws.Range("F94").Value = "mixing"
$ws.Range("F95").Value = "mixing"
$ws.Range("F96").Value = "mixing"
$ws.Range("F97").Value = "mixing"
$ws.Range("F98").Value = "safety procedures"
$ws.Range("F99").Value = "safety procedures"
$ws.Range("F100").Value = "safety procedures"
$ws.Range("F101").Value = "use restrictions"
$ws.Range("F102").Value = "use restrictions"
$ws.Range("F103").Value = "use restrictions"
$ws.Range("F118").Value = "use restrictions"
$ws.Range("F119").Value = "use restrictions"
$ws.Range("F120").Value = "application instructions"
$ws.Range("F121").Value = "application instructions"
$ws.Range("F122").Value = "use restrictions"
$ws.Range("F123").Value = "application instructions"
$ws.Range("F125").Value = "application instructions"
$ws.Range("F139").Value = "application instructions"
$ws.Range("F140").Value = "application instructions"
$ws.Range("F142").Value = "mixing"
$ws.Range("F143").Value = "mixing"
$ws.Range("F144").Value = "mixing"
$ws.Range("F145").Value = "mixing"
$ws.Range("F185").Value = "mixing"
$ws.Range("F187").Value = "mixing"
$ws.Range("F191").Value = "mixing"
$ws.Range("F197").Value = "mixing"
$ws.Range("F198").Value = "mixing"
$ws.Range("F200").Value = "mixing"
$ws.Range("F201").Value = "mixing"
$ws.Range("F202").Value = "application instructions"
$ws.Range("F203").Value = "application instructions"
$ws.Range("F204").Value = "application instructions"
$ws.Range("F205").Value = "application instructions"
$ws.Range("F206").Value = "mixing"
$ws.Range("F207").Value = "mixing"
$ws.Range("F208").Value = "mixing"
$ws.Range("F209").Value = "mixing"
$ws.Range("F210").Value = "use restrictions"
$ws.Range("F213").Value = "use restrictions"
$ws.Range("F214").Value = "use restrictions"
$ws.Range("F215").Value = "application instructions"
$ws.Range("F216").Value = "use restrictions"
$ws.Range("F217").Value = "application instructions"
$ws.Range("F218").Value = "application instructions"
$ws.Range("F220").Value = "154_pesticide_storage"
